$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.637.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "'1.843.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'315.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4267"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.61%  "
$ws.Range("D8").Value = "'0.3643"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "'0.07290"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").Value = "'0.8950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.59%  "
$ws.Range("D12").Value = "'20.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").Value = "'1.885.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'5.381"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "'6.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").Value = "'0.06889"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'78.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.58%  "
$ws.Range("D19").Value = "'0.000008874"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("D22").Value = "'27.648.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").Value = "'4.982"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "'10.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").Value = "'2.103.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'2.042"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'154.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'18.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").Value = "'1.841"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.57%  "
$ws.Range("D31").Value = "'112.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").Value = "'0.08891"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "'0.7756"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").Value = "'4.573"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.54%  "
$ws.Range("D35").Value = "'2.988"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").Value = "'1.099"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.27%  "
$ws.Range("D37").Value = "'0.9998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'0.05436"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").Value = "'1.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").Value = "'0.01928"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "'2.780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("D42").Value = "'0.5064"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").Value = "'6.805"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.19%  "
$ws.Range("D44").Value = "'0.1644"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").Value = "'8.232"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.14%  "
$ws.Range("D46").Value = "'0.06638"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "'106.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.4711"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'1.633"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.80%  "
